$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC32_Verify_store_location")

# Remove the obsolete "VERIFY_TEXT_PRESENT" step (row 24): Keyword=VERIFY_TEXT_PRESENT,
# Object=ContactUsConfirmation, ObjectType=CSS, Data_descriptor=Confirmation.
# Deleting the row shifts the subsequent rows (MyaccountSection / Logout) up by one.
$ws.Rows.Item(24).Delete()

$ws.Activate()
$ws.Rows.Item(24).Select()
